# Reorder the "UK" policy schedule sheet to reflect Daria's latest estimates,
# and drop the trailing rows (2027-2029) that are no longer part of the list.
#
# Strategy: the new row order is a reordering (subset) of the existing rows,
# so we stage each existing row's A:C values off to one side (columns AA:AC),
# then paste them back into the target rows in the new order, clear the
# staging area, and finally delete the now-unused trailing rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("UK")

$lastOldRow = 20

# Phase 1: stage current A:C values of every data row (2..20) into AA:AC
# on the same row number, preserving the shared-string (text) typing.
for ($r = 2; $r -le $lastOldRow; $r++) {
    $src = $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 3))
    $src.Copy() | Out-Null
    $dst = $ws.Range($ws.Cells.Item($r, 27), $ws.Cells.Item($r, 29))
    $dst.PasteSpecial() | Out-Null
}

# Phase 2: map new row -> old row (values taken from the staged AA:AC copy)
# and write them back into A:C for the new row order.
$rowMap = @{
    2  = 5
    3  = 13
    4  = 14
    5  = 6
    6  = 12
    7  = 8
    8  = 7
    9  = 11
    10 = 16
    11 = 4
    12 = 3
    13 = 15
    14 = 10
    15 = 17
    16 = 9
    17 = 2
}

foreach ($newRow in $rowMap.Keys) {
    $oldRow = $rowMap[$newRow]
    $src = $ws.Range($ws.Cells.Item($oldRow, 27), $ws.Cells.Item($oldRow, 29))
    $src.Copy() | Out-Null
    $dst = $ws.Range($ws.Cells.Item($newRow, 1), $ws.Cells.Item($newRow, 3))
    $dst.PasteSpecial() | Out-Null
}

# Phase 3: clear the staging area now that every value has been relocated.
$stage = $ws.Range($ws.Cells.Item(2, 27), $ws.Cells.Item($lastOldRow, 29))
$stage.Clear() | Out-Null

# Phase 4: remove the rows that are no longer part of the (now 16-row) table.
$ws.Rows.Item(18).EntireRow.Delete() | Out-Null
$ws.Rows.Item(18).EntireRow.Delete() | Out-Null
$ws.Rows.Item(18).EntireRow.Delete() | Out-Null
